$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "EventTypes"
$ws.Name = "EventTypes"

# Move the active selection on the sheet from C4 to F30
$ws.Range("F30").Select()

# Resize / reposition the workbook window to match the new saved window
# geometry (best effort - mirrors the recorded window bounds).
$win = $excel.ActiveWindow
$win.Left = -19320
$win.Top = -120
$win.Width = 19440
$win.Height = 15000
